# "añadir nombre a las restricciones y terminar el diccionario completo de lectura de datos"
# 1) sheet "warehouse": insert a new first row containing the label "warehouse"
#    (names the previously-unnamed list of locations -- the "nombre a las restricciones").
# 2) sheet "comp_quantity_inst1": flip several existing "required" flags to 1 and append
#    the remaining origin/destination combinations that were still missing, completing the
#    full dictionary (F1/F2/F3 x T1..T5, plus T5->T4) -- "terminar el diccionario completo".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# warehouse sheet: add a "warehouse" header/name above the existing T1..F3 list
# ---------------------------------------------------------------------------
$wsWarehouse = $wb.Worksheets.Item("warehouse")
$wsWarehouse.Rows.Item(1).Insert() | Out-Null
$wsWarehouse.Range("A1").Value = "warehouse"

# ---------------------------------------------------------------------------
# comp_quantity_inst1 sheet: complete the "required" column for the existing
# rows, and append the rest of the combinations so every origin/destination
# pair is represented.
# ---------------------------------------------------------------------------
$wsComp = $wb.Worksheets.Item("comp_quantity_inst1")

# Existing rows whose "required" (column D) flips from 0 to 1.
$requiredFixRows = @(3, 6, 7, 10, 11)
foreach ($r in $requiredFixRows) {
    $wsComp.Cells.Item($r, 4).Value = 1
}

# New rows 13-23 completing the dictionary.
$newRows = @(
    @(13, "F1", "T1", 18, 0),
    @(14, "F1", "T3", 102, 0),
    @(15, "F1", "T5", 105, 0),
    @(16, "F2", "T1", 60, 0),
    @(17, "F2", "T3", 58, 0),
    @(18, "F2", "T4", 91, 0),
    @(19, "F2", "T5", 20, 0),
    @(20, "F3", "T1", 20, 0),
    @(21, "F3", "T4", 34, 0),
    @(22, "F3", "T5", 61, 0),
    @(23, "T5", "T4", 2, 0)
)
foreach ($row in $newRows) {
    $r = $row[0]
    $wsComp.Cells.Item($r, 1).Value = $row[1]
    $wsComp.Cells.Item($r, 2).Value = $row[2]
    $wsComp.Cells.Item($r, 3).Value = $row[3]
    $wsComp.Cells.Item($r, 4).Value = $row[4]
}

# ---------------------------------------------------------------------------
# View/selection state, best-effort, to mirror the saved workbook UI state.
# ---------------------------------------------------------------------------
$wsParameters = $wb.Worksheets.Item("parameters")
$wsParameters.Range("E38").Select() | Out-Null

$wsWarehouse.Range("A2:A9").Select() | Out-Null

$wsComp.Range("A23:XFD23").Select() | Out-Null

$wsTrip = $wb.Worksheets.Item("trip_duration")
$wsTrip.Range("D3").Select() | Out-Null

"done"
